$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.220195055007935
$ws.Range("B1").Value = 6.182280540466309
$ws.Range("C1").Value = 2.29206371307373
$ws.Range("D1").Value = 1.501962423324585
$ws.Range("E1").Value = 1.235244154930115
